# Update Work Week and Social Spending
#
# Revises the "GDP per Capita" figures on the "Data" sheet for the years
# already present (1950-2010) and appends newly published figures for
# 2011-2016. Values are written as text (matching the source workbook,
# which stores the indicator readings as strings rather than numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$yearValues = @(
    @{ Year = 1950; Value = "1200" },
    @{ Year = 1951; Value = "1262" },
    @{ Year = 1952; Value = "1323" },
    @{ Year = 1953; Value = "1325" },
    @{ Year = 1954; Value = "1382" },
    @{ Year = 1955; Value = "1379" },
    @{ Year = 1956; Value = "1309" },
    @{ Year = 1957; Value = "1323" },
    @{ Year = 1958; Value = "1270" },
    @{ Year = 1959; Value = "1288" },
    @{ Year = 1960; Value = "1307" },
    @{ Year = 1961; Value = "1313" },
    @{ Year = 1962; Value = "1349" },
    @{ Year = 1963; Value = "1433" },
    @{ Year = 1964; Value = "1451" },
    @{ Year = 1965; Value = "1505" },
    @{ Year = 1966; Value = "1414" },
    @{ Year = 1967; Value = "1160" },
    @{ Year = 1968; Value = "1114" },
    @{ Year = 1969; Value = "1372" },
    @{ Year = 1970; Value = "1744" },
    @{ Year = 1971; Value = "1894" },
    @{ Year = 1972; Value = "1908" },
    @{ Year = 1973; Value = "2012" },
    @{ Year = 1974; Value = "2179" },
    @{ Year = 1975; Value = "2051" },
    @{ Year = 1976; Value = "2208" },
    @{ Year = 1977; Value = "2220" },
    @{ Year = 1978; Value = "2028" },
    @{ Year = 1979; Value = "2104" },
    @{ Year = 1980; Value = "2080" },
    @{ Year = 1981; Value = "1855" },
    @{ Year = 1982; Value = "1784" },
    @{ Year = 1983; Value = "1631" },
    @{ Year = 1984; Value = "1527" },
    @{ Year = 1985; Value = "1621" },
    @{ Year = 1986; Value = "1610" },
    @{ Year = 1987; Value = "1556" },
    @{ Year = 1988; Value = "1667" },
    @{ Year = 1989; Value = "1729" },
    @{ Year = 1990; Value = "1773" },
    @{ Year = 1991; Value = "1759.00349696734" },
    @{ Year = 1992; Value = "1793.94013667409" },
    @{ Year = 1993; Value = "1817.74951018463" },
    @{ Year = 1994; Value = "1820.30462461131" },
    @{ Year = 1995; Value = "1851.03748058961" },
    @{ Year = 1996; Value = "1924.08840384563" },
    @{ Year = 1997; Value = "1977.28562911474" },
    @{ Year = 1998; Value = "2025.86492178969" },
    @{ Year = 1999; Value = "2034.82752406394" },
    @{ Year = 2000; Value = "2144.82540433565" },
    @{ Year = 2001; Value = "2282.82301829164" },
    @{ Year = 2002; Value = "2609.12944509355" },
    @{ Year = 2003; Value = "2848.28129737176" },
    @{ Year = 2004; Value = "3136.62742586954" },
    @{ Year = 2005; Value = "3347.82046684666" },
    @{ Year = 2006; Value = "3564.63234245633" },
    @{ Year = 2007; Value = "3817.48352250735" },
    @{ Year = 2008; Value = "4084.20026035943" },
    @{ Year = 2009; Value = "4416.34699161145" },
    @{ Year = 2010; Value = "4904.86586262996" },
    @{ Year = 2011; Value = "5136" },
    @{ Year = 2012; Value = "5223" },
    @{ Year = 2013; Value = "5370" },
    @{ Year = 2014; Value = "5569" },
    @{ Year = 2015; Value = "5578" },
    @{ Year = 2016; Value = "5360" }
)

# Build a lookup from year -> row number for the existing rows (header is
# row 1; data originally spans rows 2-62, covering years 1950-2010).
$yearToRow = @{}
for ($r = 2; $r -le 62; $r++) {
    $y = $ws.Cells.Item($r, 4).Value()
    if ($y -ne $null) {
        $yearToRow[[int]$y] = $r
    }
}

$nextRow = 63

foreach ($entry in $yearValues) {
    $year = $entry.Year
    $value = $entry.Value

    if ($yearToRow.ContainsKey($year)) {
        $row = $yearToRow[$year]
    } else {
        $row = $nextRow
        $nextRow = $nextRow + 1

        $ws.Cells.Item($row, 1).Value = 566
        $ws.Cells.Item($row, 2).Value = "Nigeria"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }

    # Store the reading as text (leading apostrophe forces text type), then
    # strip the resulting "number stored as text" formatting so the cell
    # ends up as a plain shared-string cell, matching the source data.
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = "'" + $value
    $cell.ClearFormats()
}
